# #5: cash & deposit done
# Reshape the "存款" (deposit) sheet (sheet3): drop the old exchange-rate
# "quantity" column, move the total amount left into column F, and append
# the standard property/legislator metadata columns (G:M) that the other
# sheets already carry. Row 1 becomes a real header row of field names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# ---- header row (row 1) ---------------------------------------------
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# ---- data rows (2-12) -------------------------------------------------
# column A (index) stays as-is; B/C/D/E stay as-is; the former F (quantity)
# is dropped, the former G (total) shifts left into F, and G..M are new
# metadata columns shared by the rest of the workbook's sheets.

$rows = 75,76,77,78,79,80,81,82,83,84,85
$totals = @{
    75 = 5474321
    76 = 1000000
    77 = 160800
    78 = 401159
    79 = 1306261
    80 = 3000000
    81 = 182
    82 = 369
    83 = 1194111
    84 = 769060
    85 = 199408
}

for ($r = 2; $r -le 12; $r++) {
    $idx = $rows[$r - 2]

    # drop the old exchange-rate "quantity" value, replace with the total
    $ws.Cells.Item($r, 6).Value = $totals[$idx]

    $ws.Cells.Item($r, 7).Value = "deposit"
    $ws.Cells.Item($r, 8).Value = "normal"
    $ws.Cells.Item($r, 9).Value = "2012-04-16"
    $ws.Cells.Item($r, 10).Value = "呂學樟"
    $ws.Cells.Item($r, 11).Value = 892
    $ws.Cells.Item($r, 12).Value = "tmp63271"
    $ws.Cells.Item($r, 13).Value = $idx
}
